$d = $word.ActiveDocument

# Locate the paragraph that starts the block to be replaced ("<tab>Indent 1.1.")
# and the paragraph that starts right after the block ("<tab>Indent 1.2."),
# then overwrite everything in between (the two "Indent level 2.x" lines and
# the two "Numbered item under plain text" lines) with the new tabbed
# hierarchy that introduces the "Fully Bold Header" sections.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if (($null -eq $startPara) -and ($p.Range.Text -like "*Indent 1.1.*")) {
        $startPara = $p
    }
    elseif (($null -ne $startPara) -and ($null -eq $endPara) -and ($p.Range.Text -like "*Indent 1.2.*")) {
        $endPara = $p
        break
    }
}

if ($null -eq $startPara -or $null -eq $endPara) {
    throw "Could not locate the anchor paragraphs (Indent 1.1. / Indent 1.2.)"
}

$rangeStart = $startPara.Range.Start
$rangeEnd = $endPara.Range.Start
$r = $d.Range($rangeStart, $rangeEnd)

$newBlock = @'
<w:p><w:r><w:tab/></w:r><w:r><w:t>Indent 1.1.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/><w:t>Fully Bold Header</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>Text under fully bold header.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Text under &#8220;text under fuly&#8221;.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Text under &#8220; text under&#8221; 2.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Fully Bold Header</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> 2: </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>A.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>b.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:b/></w:rPr><w:tab/></w:r><w:r><w:t>c.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Text under fully.</w:t></w:r></w:p><w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Text under fully 2.</w:t></w:r></w:p>
'@

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $newBlock + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

Write-Output "Replaced tabbed-indent block with fully-bold-header hierarchy."
